$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dic_Disagg_Ausprägungen")

# Row 84 - A_SERIES_33ADD
$ws.Range("C84").Value = "Beitrag Deutschlands zur globalen Pandemieprävention und -reaktion"
$ws.Range("D84").Value = "XXXBeitrag Deutschlands zur globalen Pandemieprävention und -reaktion"

# Row 91 - A_SERIES_BEREIN
$ws.Range("C91").Value = "Bereinigter Gender Pay Gap"

# Row 108 - A_SERIES_FATFTEC
$ws.Range("C108").Value = "Technical Compliance"

# Row 136 - A_SERIES_RMC
$ws.Range("C136").Value = "Rohstoffeinsatz für Konsum und Investitionen (RMC)"
$ws.Range("D136").Value = "XXXRohstoffeinsatz für Konsum und Investitionen (RMC)"

# Row 137 - A_SERIES_RMI
$ws.Range("C137").Value = "Rohstoffeinsatz für Exporte (RMC)"
$ws.Range("D137").Value = "XXXRohstoffeinsatz für Exporte (RMC)"

# Row 146 - A_SERIES_VAETERMONAT
$ws.Range("C146").Value = "Männliche Bezieher an allen genommenen Elterngeldmonaten"
